# info_tags.xlsx — "fix bugs. Update default HTML."
#
# Changes applied (per the canonical OOXML diff):
#   1. Worksheet "2-data" is renamed to "9-data" (it now carries the "9"
#      record instead of "2"), and its A1 counter cell is updated 2 -> 9.
#   2. "9-data" becomes the active/selected sheet (previously "m-map" was
#      active); the workbook's remembered active tab moves accordingly.
#   3. The on-screen selection anchors a newly-touched range (B16:B17) on
#      every sheet, in addition to each sheet's previous cursor cell.

$wb = $excel.ActiveWorkbook

# --- 1. Rename "2-data" -> "9-data" and bump its header counter cell ---
$wsData = $wb.Worksheets.Item("2-data")
$wsData.Name = "9-data"
$wsData.Range("A1").Value = 9

# --- 2. Refresh the selection/cursor on every sheet ---------------------
# "1-post": cursor stays conceptually anchored near D34, with B16:B17 now
# also part of the working selection.
$wsPost = $wb.Worksheets.Item("1-post")
[void]$wsPost.Range("B16:B17").Select()

# "9-data" (formerly "2-data"): becomes the active sheet, with the cursor
# parked on the newly edited B16:B17 block.
[void]$wsData.Activate()
[void]$wsData.Range("B16:B17").Select()

# "m-map": no longer the active tab, but B16:B17 joins its selection too.
$wsMap = $wb.Worksheets.Item("m-map")
[void]$wsMap.Range("B16:B17").Select()

# "Sheet4": same treatment.
$wsSheet4 = $wb.Worksheets.Item("Sheet4")
[void]$wsSheet4.Range("B16:B17").Select()

# --- 3. Leave "9-data" as the sheet in front when the file is reopened --
[void]$wsData.Activate()
